$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032960743110995
$ws.Range("D2").Value = 1.049893678851363
$ws.Range("E2").Value = 1.043039596194371
$ws.Range("F2").Value = 1.055671848741947
$ws.Range("I2").Value = 1.042264514150976
$ws.Range("J2").Value = 1.038088090622999
$ws.Range("K2").Value = 1.052649469556411
$ws.Range("L2").Value = 1.045814585829269
$ws.Range("M2").Value = 1.058411664113584
$ws.Range("N2").Value = 1.016678052744064

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033747246808538
$ws.Range("D3").Value = 1.050477735552251
$ws.Range("E3").Value = 1.043724726960093
$ws.Range("F3").Value = 1.05640656365925
$ws.Range("I3").Value = 1.042441859091367
$ws.Range("J3").Value = 1.038517935411993
$ws.Range("K3").Value = 1.053046452851838
$ws.Range("L3").Value = 1.046310999899279
$ws.Range("M3").Value = 1.058960067476547
$ws.Range("N3").Value = 1.016821079793987

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.034256856724133
$ws.Range("D4").Value = 1.050856011286234
$ws.Range("E4").Value = 1.044169021881312
$ws.Range("F4").Value = 1.056882902892485
$ws.Range("I4").Value = 1.042555507434947
$ws.Range("J4").Value = 1.038796080100585
$ws.Range("K4").Value = 1.053302978216569
$ws.Range("L4").Value = 1.046632501993623
$ws.Range("M4").Value = 1.059315179211594
$ws.Range("N4").Value = 1.016913607024729

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034471259905748
$ws.Range("D5").Value = 1.051015120614056
$ws.Range("E5").Value = 1.044356033535002
$ws.Range("F5").Value = 1.05708337623403
$ws.Range("I5").Value = 1.042603019563173
$ws.Range("J5").Value = 1.038913012301895
$ws.Range("K5").Value = 1.05341073630811
$ws.Range("L5").Value = 1.046767729102488
$ws.Range("M5").Value = 1.059464528183187
$ws.Range("N5").Value = 1.016952500003578

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034507268657329
$ws.Range("D6").Value = 1.051041840531142
$ws.Range("E6").Value = 1.04438744704674
$ws.Range("F6").Value = 1.057117049446034
$ws.Range("I6").Value = 1.042610981454522
$ws.Range("J6").Value = 1.03893264570396
$ws.Range("K6").Value = 1.053428824326751
$ws.Range("L6").Value = 1.046790438249513
$ws.Range("M6").Value = 1.059489607998942
$ws.Range("N6").Value = 1.01695902997116

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.034259720951134
$ws.Range("D7").Value = 1.05085813699243
$ws.Range("E7").Value = 1.044171519838729
$ws.Range("F7").Value = 1.056885580764131
$ws.Range("I7").Value = 1.04255614333926
$ws.Range("J7").Value = 1.038797642555394
$ws.Range("K7").Value = 1.053304418422056
$ws.Range("L7").Value = 1.046634308641098
$ws.Range("M7").Value = 1.059317174586362
$ws.Range("N7").Value = 1.01691412673671

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03322640183709
$ws.Range("D8").Value = 1.050090989428937
$ws.Range("E8").Value = 1.043270937381792
$ws.Range("F8").Value = 1.055919955291077
$ws.Range("I8").Value = 1.042324677085119
$ws.Range("J8").Value = 1.038233356552206
$ws.Range("K8").Value = 1.052783703339862
$ws.Range("L8").Value = 1.045982290513611
$ws.Range("M8").Value = 1.058596945084107
$ws.Range("N8").Value = 1.016726393430077

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031410924882021
$ws.Range("D9").Value = 1.048741964175743
$ws.Range("E9").Value = 1.041691511055724
$ws.Range("F9").Value = 1.054225616714031
$ws.Range("I9").Value = 1.04190837685561
$ws.Range("J9").Value = 1.037239120965103
$ws.Range("K9").Value = 1.051863528702778
$ws.Range("L9").Value = 1.04483563349354
$ws.Range("M9").Value = 1.057329864374096
$ws.Range("N9").Value = 1.016395445049912

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030204325133539
$ws.Range("D10").Value = 1.047844610616016
$ws.Range("E10").Value = 1.040643734295759
$ws.Range("F10").Value = 1.053101041899885
$ws.Range("I10").Value = 1.041625233042663
$ws.Range("J10").Value = 1.036576450776833
$ws.Range("K10").Value = 1.051248416809834
$ws.Range("L10").Value = 1.044072822990617
$ws.Range("M10").Value = 1.056486632889325
$ws.Range("N10").Value = 1.016174748527448

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029682759745985
$ws.Range("D11").Value = 1.047456546247752
$ws.Range("E11").Value = 1.040191288376832
$ws.Range("F11").Value = 1.052615298774534
$ws.Range("I11").Value = 1.041501309692325
$ws.Range("J11").Value = 1.036289560742214
$ws.Range("K11").Value = 1.050981691687663
$ws.Range("L11").Value = 1.043742923357116
$ws.Range("M11").Value = 1.056121880151877
$ws.Range("N11").Value = 1.016079175330187

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029489164311208
$ws.Range("D12").Value = 1.04731247844866
$ws.Range("E12").Value = 1.040023419467094
$ws.Range("F12").Value = 1.052435055498815
$ws.Range("I12").Value = 1.041455081629945
$ws.Range("J12").Value = 1.03618300598035
$ws.Range("K12").Value = 1.050882562933505
$ws.Range("L12").Value = 1.043620445993293
$ws.Range("M12").Value = 1.055986452480572
$ws.Range("N12").Value = 1.01604367415652

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029530684950973
$ws.Range("D13").Value = 1.04734337798095
$ws.Range("E13").Value = 1.040059419302325
$ws.Range("F13").Value = 1.052473709959787
$ws.Range("I13").Value = 1.041465006624573
$ws.Range("J13").Value = 1.036205861910402
$ws.Range("K13").Value = 1.050903828859147
$ws.Range("L13").Value = 1.043646714970079
$ws.Range("M13").Value = 1.056015499533631
$ws.Range("N13").Value = 1.016051289317635

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029666754277572
$ws.Range("D14").Value = 1.047444635983965
$ws.Range("E14").Value = 1.040177408399449
$ws.Range("F14").Value = 1.052600396055694
$ws.Range("I14").Value = 1.041497492488016
$ws.Range("J14").Value = 1.036280752701386
$ws.Range("K14").Value = 1.050973498790292
$ws.Range("L14").Value = 1.043732798064991
$ws.Range("M14").Value = 1.056110684469905
$ws.Range("N14").Value = 1.016076240810326

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029750609314674
$ws.Range("D15").Value = 1.047507034593264
$ws.Range("E15").Value = 1.040250130569542
$ws.Range("F15").Value = 1.052678475904236
$ws.Range("I15").Value = 1.041517481969548
$ws.Range("J15").Value = 1.036326896619364
$ws.Range("K15").Value = 1.051016417470733
$ws.Range("L15").Value = 1.043785844970771
$ws.Range("M15").Value = 1.056169338756149
$ws.Range("N15").Value = 1.01609161412631

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03023895884184
$ws.Range("D16").Value = 1.047870375793902
$ws.Range("E16").Value = 1.040673788121741
$ws.Range("F16").Value = 1.053133304658265
$ws.Range("I16").Value = 1.04163342969172
$ws.Range("J16").Value = 1.036595491905947
$ws.Range("K16").Value = 1.051266110655466
$ws.Range("L16").Value = 1.044094725956487
$ws.Range("M16").Value = 1.056510848323079
$ws.Range("N16").Value = 1.01618109123399

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030545529978572
$ws.Range("D17").Value = 1.048098424326771
$ws.Range("E17").Value = 1.040939873010456
$ws.Range("F17").Value = 1.05341893099555
$ws.Range("I17").Value = 1.041705807836144
$ws.Range("J17").Value = 1.036763989135365
$ws.Range("K17").Value = 1.051422636600775
$ws.Range("L17").Value = 1.044288587688978
$ws.Range("M17").Value = 1.056725169161483
$ws.Range("N17").Value = 1.016237215490421

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030724434469759
$ws.Range("D18").Value = 1.048231488858317
$ws.Range("E18").Value = 1.041095196129485
$ws.Range("F18").Value = 1.053585648136212
$ws.Range("I18").Value = 1.041747897339369
$ws.Range("J18").Value = 1.036862275501515
$ws.Range("K18").Value = 1.051513898997492
$ws.Range("L18").Value = 1.044401702671205
$ws.Range("M18").Value = 1.056850214655786
$ws.Range("N18").Value = 1.016269950807717

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030785450931185
$ws.Range("D19").Value = 1.048276868447732
$ws.Range("E19").Value = 1.041148177641015
$ws.Range("F19").Value = 1.05364251400107
$ws.Range("I19").Value = 1.041762227127161
$ws.Range("J19").Value = 1.03689578938172
$ws.Range("K19").Value = 1.051545010871066
$ws.Range("L19").Value = 1.044440278474748
$ws.Range("M19").Value = 1.056892857950017
$ws.Range("N19").Value = 1.016281112518349

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030512628811835
$ws.Range("D20").Value = 1.048073951930298
$ws.Range("E20").Value = 1.040911312161205
$ws.Range("F20").Value = 1.053388273964312
$ws.Range("I20").Value = 1.041698055521608
$ws.Range("J20").Value = 1.036745910472617
$ws.Range("K20").Value = 1.051405846616304
$ws.Range("L20").Value = 1.044267784143903
$ws.Range("M20").Value = 1.056702170850172
$ws.Range("N20").Value = 1.016231193989983

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029626681450275
$ws.Range("D21").Value = 1.047414815894719
$ws.Range("E21").Value = 1.040142658299028
$ws.Range("F21").Value = 1.052563085078385
$ws.Range("I21").Value = 1.041487931652965
$ws.Range("J21").Value = 1.036258698971338
$ws.Range("K21").Value = 1.050952984242403
$ws.Range("L21").Value = 1.043707447013411
$ws.Range("M21").Value = 1.056082653270052
$ws.Range("N21").Value = 1.016068893242194

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02907044577119
$ws.Range("D22").Value = 1.047000835055967
$ws.Range("E22").Value = 1.039660473400191
$ws.Range("F22").Value = 1.052045318000175
$ws.Range("I22").Value = 1.041354676732488
$ws.Range("J22").Value = 1.035952422346978
$ws.Range("K22").Value = 1.05066793328199
$ws.Range("L22").Value = 1.043355500605696
$ws.Range("M22").Value = 1.055693472929137
$ws.Range("N22").Value = 1.01596684258765

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029365240870194
$ws.Range("D23").Value = 1.047220251324566
$ws.Range("E23").Value = 1.039915983932153
$ws.Range("F23").Value = 1.052319694705453
$ws.Range("I23").Value = 1.041425425595434
$ws.Range("J23").Value = 1.03611477993596
$ws.Range("K23").Value = 1.050819073849615
$ws.Range("L23").Value = 1.043542039418112
$ws.Range("M23").Value = 1.055899752429003
$ws.Range("N23").Value = 1.016020941949926

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030527495165124
$ws.Range("D24").Value = 1.048085009805822
$ws.Range("E24").Value = 1.040924217208155
$ws.Range("F24").Value = 1.053402126198595
$ws.Range("I24").Value = 1.041701558852752
$ws.Range("J24").Value = 1.036754079427528
$ws.Range("K24").Value = 1.051413433400675
$ws.Range("L24").Value = 1.044277184251039
$ws.Range("M24").Value = 1.056712562686807
$ws.Range("N24").Value = 1.016233914850029

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031879622050186
$ws.Range("D25").Value = 1.049090376606059
$ws.Range("E25").Value = 1.042098927522726
$ws.Range("F25").Value = 1.054662774994585
$ws.Range("I25").Value = 1.042016993373013
$ws.Range("J25").Value = 1.037496133532618
$ws.Range("K25").Value = 1.052101715549798
$ws.Range("L25").Value = 1.04513179176454
$ws.Range("M25").Value = 1.057657180477952
$ws.Range("N25").Value = 1.016481016567111
